$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Date" column (H) to the imported data, formatted as a short date.
# Apply the date number format to H1 first, then propagate the exact same
# style to H2:H3 via copy/paste-format so all three cells share one style
# index (matches how Excel itself would fill a column format down).
$ws.Range("H1").NumberFormat = "mm-dd-yy"
$ws.Range("H1").Copy()
$ws.Range("H2:H3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H1").Value = 41733
$ws.Range("H2").Value = 42129

# Move the active selection to reflect where the user ended up after editing.
[void]$ws.Range("H10").Select()
